# Correct the null (empty) values in the "Security" column (E) of the
# df_toyota sheet by filling them with "Information not available".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "E2","E3","E4","E5","E6","E7","E8","E9","E10","E11","E12","E13","E14",
    "E15","E16","E17","E18","E19","E20","E21","E22","E23","E27","E28","E29",
    "E30","E31","E37","E38","E39","E40","E41","E42","E50","E51","E52","E53",
    "E54","E66","E67","E68","E69","E70","E82","E83","E84","E85","E86","E87",
    "E99","E100","E101","E102","E103","E104","E116","E117","E118","E134",
    "E135","E150","E151"
)

foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    if ([string]::IsNullOrEmpty($cell.Value2)) {
        $cell.Value2 = "Information not available"
    }
}
